$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = "Chevrolet Corsa"
$ws.Range("C4").Value = "Corsa Classic"
$ws.Range("D4").Value = "https://s2.glbimg.com/nFF9NwKwQwBb6LmPpK7bPBzm4ho=/620x465/s.glbimg.com/jo/g1/f/original/2014/06/04/classic_2015.jpg"
$ws.Range("E4").Value = "PH4701"
$ws.Range("F4").Value = "CA5496"
$ws.Range("G4").Value = "G5995"
$ws.Range("H4").Value = "HK301"
$ws.Range("I4").Value = "Chevrolet"
$ws.Range("J4").Value = "Chevrolet"
$ws.Range("L4").Value = "Agile/Astra/Celta/Classic/Corsa"
$ws.Range("L1").Value = "detail"

# Scroll the sheet so column D is the leftmost visible column, then land the
# cursor on E11 (matches the saved view state in the target workbook).
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("E11").Select() | Out-Null
